# chore: add monthly employment outputs
#
# Refresh pass over the employment_master workbook:
#   - every data row (2..50) on every sheet gets its collected_at (col J)
#     timestamp bumped to the latest collection run
#   - the 피보험자수 (insured-count) sheet's first 9 regions (rows 2..10)
#     get their current_value (col E) updated with the newly collected
#     figures; 광명시 (row 5) also flips from 정상 -> 주의 on all three
#     signal columns (F current_signal, G prev_1m_signal, I prev_2m_signal)

$wb = $excel.ActiveWorkbook

$newTimestamp = "2026-02-12T23:29:38"
$firstDataRow = 2
$lastDataRow = 50
$timestampCol = 10   # J

# 1) Bump collected_at on every row of every sheet.
for ($s = 1; $s -le $wb.Worksheets.Count; $s++) {
    $ws = $wb.Worksheets.Item($s)
    for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
        $ws.Cells.Item($r, $timestampCol).Value = $newTimestamp
    }
}

# 2) Updated current_value figures for 피보험자수 (insured count), rows 2-10.
$insured = $wb.Worksheets.Item(6)

$currentValueCol = 5  # E

$newValues = @{
    2  = 13589
    3  = 188692
    4  = 56486
    5  = 68592
    6  = 91189
    7  = 32121
    8  = 59443
    9  = 129118
    10 = 120613
}

foreach ($row in $newValues.Keys) {
    $insured.Cells.Item($row, $currentValueCol).Value = $newValues[$row]
}

# 3) 광명시 (row 5) moves into the "주의" (caution) signal band.
$insured.Cells.Item(5, 6).Value = "주의"   # F5 current_signal
$insured.Cells.Item(5, 7).Value = "주의"   # G5 prev_1m_signal
$insured.Cells.Item(5, 9).Value = "주의"   # I5 prev_2m_signal
